$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (B2:D2 keep same text, but some numeric columns change)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.703265666666667
$ws.Range("H2").Value = 5.109797
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 0.1642713333333333
$ws.Range("Q2").Value = 0.2797977220842222
$ws.Range("R2").Value = 2.518179498758
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

# Delete row 3 entirely (the MuSCs row)
$ws.Rows("3").Delete()
